$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$zoomUrl = "https://nih.zoomgov.com/j/1618272266?pwd=U281WDVxLzcyN0VNdkd3d3kzQ2d3Zz09"

# Row 3: Intro to R: Data wrangling - 7/8/2024
$ws.Range("E3").Value = 1720443600
$ws.Range("G3").Value = $zoomUrl

# Row 4: Intro to R: Data visualization - 7/15/2024
$ws.Range("E4").Value = 1721048400
$ws.Range("G4").Value = $zoomUrl

# Row 5: Intro to R: Data analysis - 7/22/2024
$ws.Range("E5").Value = 1721653200
$ws.Range("G5").Value = $zoomUrl

# Row 6: Real-world data analysis in R - 7/29/2024
$ws.Range("E6").Value = 1722258000
$ws.Range("G6").Value = $zoomUrl

# Remove the old "Statistical hypothesis testing in Prism" row (row 7) entirely
$ws.Range("A7").EntireRow.Delete()

# Update sheet view: drop the leftmost-visible-column hint, select the new blank row
[void]$ws.Range("A1").Select()
[void]$ws.Range("A7:XFD7").Select()
